# "arrival read successful, with data cleaned"
# The arrival-time series in column A was supposed to be sorted ascending;
# row 24 held an out-of-order timestamp (11:25:30 AM, serial 0.4760416666666667)
# sitting between 11:22:40 AM (A23) and 11:22:57 AM (A25). Clean it up by
# replacing it with the correct reading of 11:22:30 AM so the column sorts
# correctly again.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 0.47395833333333331

# Move the live selection/active cell to where the analyst ended up after
# verifying the fix.
$ws.Range("C26").Select()
